$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns stay formatted as Text so that
# numeric-looking strings (e.g. "17.90", "0.06316", "1.001") are not
# silently coerced into floating point numbers and lose their exact
# textual representation (trailing zeros, scientific notation, etc).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.419.75"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.853.86"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "233.83"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4683"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "0.2739"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "0.06316"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "1.862.68"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "16.96"
$ws.Range("E11").Value = "  +5.23%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07458"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "84.24"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "0.6234"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "30.376.16"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "228.67"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000007311"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "4.942"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").Value = "5.907"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "167.37"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "9.195"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "17.90"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "1.881"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").Value = "0.1021"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "1.374"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "4.097"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "3.822"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").Value = "0.04892"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "1.141"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "0.7075"
$ws.Range("E34").Value = "  -3.41%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "0.01929"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").Value = "2.670"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "0.8735"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").Value = "1.950"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").Value = "105.57"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "0.9995"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "5.522"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Value = "0.4058"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "7.110"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "61.33"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "0.1215"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("B47").Value = "Elrond"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D47").Value = "33.46"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.569"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "0.05548"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").Value = "1.356"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "0.3676"
